$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 4644.1304
$ws.Range("I98").Value = 4329.5835
$ws.Range("J98").Value = 5776.5
$ws.Range("K98").Value = 4329.5835
$ws.Range("L98").Value = 5776.5
$ws.Range("M98").Value = -2831.5835
$ws.Range("N98").Value = -8772.5
$ws.Range("H107").Value = 19739592
$ws.Range("I107").Value = 8336696.5
$ws.Range("J107").Value = 62500450
$ws.Range("K107").Value = 8336696.5
$ws.Range("L107").Value = 62500450
$ws.Range("M107").Value = -8334776.5
$ws.Range("N107").Value = -62504290
$ws.Range("H116").Value = 27783866
$ws.Range("I116").Value = 125000500
$ws.Range("J116").Value = 7685.7144
$ws.Range("K116").Value = 125000500
$ws.Range("L116").Value = 7685.7144
$ws.Range("M116").Value = -124997058
$ws.Range("N116").Value = -14569.7144
$ws.Range("H122").Value = 4644.1304
$ws.Range("I122").Value = 4329.5835
$ws.Range("J122").Value = 5776.5
$ws.Range("K122").Value = 12988.7505
$ws.Range("L122").Value = 17329.5
$ws.Range("M122").Value = -10538.7505
$ws.Range("N122").Value = -22229.5
$ws.Range("H137").Value = 9793.5
$ws.Range("I137").Value = 13988.4
$ws.Range("J137").Value = 5598.6
$ws.Range("K137").Value = 41965.2
$ws.Range("L137").Value = 16795.8
$ws.Range("M137").Value = -39415.2
$ws.Range("N137").Value = -21895.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4527.606
$ws.Range("I61").Value = 3169.8462
$ws.Range("J61").Value = 9570.714
$ws.Range("K61").Value = 3169.8462
$ws.Range("L61").Value = 9570.714
$ws.Range("M61").Value = -2957.8462
$ws.Range("N61").Value = -9994.714
$ws.Range("H74").Value = 166901.2
$ws.Range("I74").Value = 269499.66
$ws.Range("J74").Value = 13003.5
$ws.Range("K74").Value = 269499.66
$ws.Range("L74").Value = 13003.5
$ws.Range("M74").Value = -268625.66
$ws.Range("N74").Value = -14751.5
$ws.Range("H77").Value = 166901.2
$ws.Range("I77").Value = 269499.66
$ws.Range("J77").Value = 13003.5
$ws.Range("K77").Value = 1347498.3
$ws.Range("L77").Value = 65017.5
$ws.Range("M77").Value = -1343130.3
$ws.Range("N77").Value = -73753.5
$ws.Range("H136").Value = 4527.606
$ws.Range("I136").Value = 3169.8462
$ws.Range("J136").Value = 9570.714
$ws.Range("K136").Value = 9509.5386
$ws.Range("L136").Value = 28712.142
$ws.Range("M136").Value = -6959.5386
$ws.Range("N136").Value = -33812.142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 6483.5
$ws.Range("I22").Value = 6907.7334
$ws.Range("J22").Value = 120
$ws.Range("K22").Value = 6907.7334
$ws.Range("L22").Value = 120
$ws.Range("M22").Value = -6734.7334
$ws.Range("N22").Value = -466

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15483.182
$ws.Range("I31").Value = 6865.2
$ws.Range("J31").Value = 22664.834
$ws.Range("K31").Value = 6865.2
$ws.Range("L31").Value = 22664.834
$ws.Range("M31").Value = -6570.2
$ws.Range("N31").Value = -23254.834
$ws.Range("H34").Value = 15483.182
$ws.Range("I34").Value = 6865.2
$ws.Range("J34").Value = 22664.834
$ws.Range("K34").Value = 6865.2
$ws.Range("L34").Value = 22664.834
$ws.Range("M34").Value = -6663.2
$ws.Range("N34").Value = -23068.834
$ws.Range("H99").Value = 6902.1665
$ws.Range("J99").Value = 6736.7144
$ws.Range("L99").Value = 6736.7144
$ws.Range("N99").Value = -9732.714400000001
$ws.Range("H126").Value = 6902.1665
$ws.Range("J126").Value = 6736.7144
$ws.Range("L126").Value = 20210.1432
$ws.Range("N126").Value = -25150.1432
$ws.Range("H141").Value = 95000
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 8001997.5
$ws.Range("J5").Value = 1995.6666
$ws.Range("L5").Value = 5986.9998
$ws.Range("N5").Value = -6210.9998
$ws.Range("H107").Value = 16667516
$ws.Range("J107").Value = 25000950
$ws.Range("L107").Value = 75002850
$ws.Range("N107").Value = -75006690
$ws.Range("H122").Value = 1770052.2
$ws.Range("I122").Value = 3144659
$ws.Range("J122").Value = 2700.5715
$ws.Range("K122").Value = 28301931
$ws.Range("L122").Value = 24305.1435
$ws.Range("M122").Value = -28299481
$ws.Range("N122").Value = -29205.1435
$ws.Range("H135").Value = 8001997.5
$ws.Range("J135").Value = 1995.6666
$ws.Range("L135").Value = 17960.9994
$ws.Range("N135").Value = -23030.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 651.5714
$ws.Range("I102").Value = 547.8461
$ws.Range("K102").Value = 547.8461
$ws.Range("M102").Value = 1074.1539
$ws.Range("H113").Value = 5646.826
$ws.Range("I113").Value = 3687.3044
$ws.Range("K113").Value = 3687.3044
$ws.Range("M113").Value = -1517.3044

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5195.905
$ws.Range("I7").Value = 3300.9
$ws.Range("K7").Value = 3300.9
$ws.Range("M7").Value = -3188.9
$ws.Range("H40").Value = 4982.968
$ws.Range("I40").Value = 4261.6816
$ws.Range("K40").Value = 4261.6816
$ws.Range("M40").Value = -4125.6816
$ws.Range("H61").Value = 3201.1562
$ws.Range("I61").Value = 2294.6667
$ws.Range("J61").Value = 4931.727
$ws.Range("K61").Value = 2294.6667
$ws.Range("L61").Value = 4931.727
$ws.Range("M61").Value = -2092.6667
$ws.Range("N61").Value = -5335.727
$ws.Range("H113").Value = 3201.1562
$ws.Range("I113").Value = 2294.6667
$ws.Range("J113").Value = 4931.727
$ws.Range("K113").Value = 2294.6667
$ws.Range("L113").Value = 4931.727
$ws.Range("M113").Value = -124.6667000000002
$ws.Range("N113").Value = -9271.726999999999
$ws.Range("H126").Value = 5195.905
$ws.Range("I126").Value = 3300.9
$ws.Range("K126").Value = 9902.700000000001
$ws.Range("M126").Value = -7432.700000000001
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H136").Value = 5685.7144
$ws.Range("I136").Value = 3000
$ws.Range("J136").Value = 7700
$ws.Range("K136").Value = 9000
$ws.Range("L136").Value = 23100
$ws.Range("M136").Value = -6450
$ws.Range("N136").Value = -28200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 50015452
$ws.Range("I81").Value = 900.5
$ws.Range("K81").Value = 1801
$ws.Range("M81").Value = -740
$ws.Range("H84").Value = 50015452
$ws.Range("I84").Value = 900.5
$ws.Range("K84").Value = 9005
$ws.Range("M84").Value = -3701
$ws.Range("H122").Value = 25207448
$ws.Range("I122").Value = 33607064
$ws.Range("K122").Value = 100821192
$ws.Range("M122").Value = -100818742
$ws.Range("H126").Value = 2304.3
$ws.Range("I126").Value = 1492.5652
$ws.Range("K126").Value = 4477.6956
$ws.Range("M126").Value = -2007.6956
$ws.Range("H136").Value = 45504036
$ws.Range("I136").Value = 71430520
$ws.Range("K136").Value = 214291560
$ws.Range("M136").Value = -214289010
